$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R ("backup"), mirroring the existing header style ---
$ws.Range("R1").Value = "backup"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# R2:R204 default to 0, with a couple of overrides
$rOverrides = @{ 59 = 1; 64 = 1 }
for ($r = 2; $r -le 204; $r++) {
    $val = 0
    if ($rOverrides.ContainsKey($r)) { $val = $rOverrides[$r] }
    $ws.Cells.Item($r, 18).Value = $val
}

# --- detect_structure (Q) rolling re-evaluation resets ---
$qResetRows = @(32, 40, 44, 54)
foreach ($r in $qResetRows) {
    $ws.Cells.Item($r, 17).Value = 0
}

# --- isPivot (O) updates for the now-not-last-bar rows ---
$ws.Cells.Item(202, 15).Value = 1
$ws.Cells.Item(204, 15).Value = 2

# --- Three new monthly bars appended ---
$newRows = @{
    205 = @(45474, 826, 891.5499877929688, 778.2000122070312, 889.1500244140625, 884.0974731445312, 84319472, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0)
    206 = @(45505, 890, 901, 800.4000244140625, 845.0999755859375, 845.0999755859375, 65001774, 2024, 8, 1, 0, 0, 0, 31, 0, 0, 2)
    207 = @(45536, 850, 929, 809.5499877929688, 914.0499877929688, 914.0499877929688, 59808274, 2024, 9, 1, 0, 0, 0, 35, 0, 0, 0)
}

foreach ($r in 205..207) {
    $vals = $newRows[$r]
    for ($c = 1; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Output "done"